$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 98 (pushes existing rows 98:151 down to 99:152)
$ws.Rows.Item(98).EntireRow.Insert()

# Populate the newly inserted row 98 with the new weekly price record
$ws.Cells.Item(98, 1).Value = 11
$ws.Cells.Item(98, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(98, 3).Value = "Bíobío"
$ws.Cells.Item(98, 4).Value = 44523
$ws.Cells.Item(98, 5).Value = 8
$ws.Cells.Item(98, 6).Value = 100114013
$ws.Cells.Item(98, 7).Value = "Zanahoria"
$ws.Cells.Item(98, 8).Value = "Sin especificar"
$ws.Cells.Item(98, 9).Value = "Primera"
$ws.Cells.Item(98, 10).Value = 1000
$ws.Cells.Item(98, 11).Value = 7000
$ws.Cells.Item(98, 12).Value = 7500
$ws.Cells.Item(98, 13).Value = 7250
$ws.Cells.Item(98, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(98, 15).Value = "Región de Ñuble"
$ws.Cells.Item(98, 16).Value = 362
$ws.Cells.Item(98, 17).Value = 20
$ws.Cells.Item(98, 18).Value = "Hortaliza"
